# Build Release-: V 4.24.0.160
#
# The roster's QCO shared-string block shifted: a new hire ("Danielle
# Gaylord") was added to the list, bumping most of the subsequent rows'
# referenced text, and the trailing "OFF" row became blank.
#
# Values are staged in a scratch row (far outside the sheet's used range)
# and then copied into place with Range.Copy — a plain `.Value = "...`n..."`
# assignment triggers the engine's auto row-height recalculation for
# multi-line text, which the target workbook does not have. Copy/paste
# does not trigger that recalculation, and the scratch row is deleted
# afterwards so it leaves no trace.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee")

# Row 9 ("OFF") becomes blank; copy the already-blank A10 cell so the
# shared-string-typed empty cell is preserved instead of the cell being
# dropped (which a plain `.Value = ""` assignment would do). Done first,
# before the scratch-row copies below, so this direct cell-to-cell copy
# is not affected by the scratch row's use of the clipboard.
$ws.Cells.Item(10, 1).Copy($ws.Cells.Item(9, 1))

$scratchRow = 500
$scratch = $ws.Cells.Item($scratchRow, 1)

function Set-RosterCell([int]$row, [string]$text) {
    $scratch.Value = $text
    $scratch.Copy($ws.Cells.Item($row, 1))
}

Set-RosterCell 2 "89212114 - Mohammed Turner`nROLE : RTGO100 2023-11-08T12:47:56.947450800"
Set-RosterCell 3 "92457737 - Lilliana Williamson`nROLE : RTGO100 1701844270281"
Set-RosterCell 4 "90317880 - Lewis Mosciski`nROLE : RTGO100 1701844270281"
Set-RosterCell 5 "90833312 - Angelo Mueller`nROLE : RTGO100 1701844270281"
Set-RosterCell 6 "92970163 - Glenna Lynch`nROLE : RTGO100 1701853905917"
Set-RosterCell 7 "68306525 - Danielle Gaylord`nROLE : QCO 2023-11-24T09:25:13.428483500"

# Remove the scratch row so it leaves no residue in the saved sheet.
$scratch.EntireRow.Delete()
